$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data (rows for "Primera" and "Segunda" quality) was added at
# the top of this block of records. Insert two rows before the existing
# row 493, shifting the existing rows 493:522 down to 495:524.
$ws.Rows("493:494").Insert()

# Fill in the two newly inserted rows with the new record data.

# Row 493 - Primera
$ws.Cells.Item(493, 1).Value = 1
$ws.Cells.Item(493, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(493, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(493, 4).Value = 45021
$ws.Cells.Item(493, 5).Value = 15
$ws.Cells.Item(493, 6).Value = 100112032
$ws.Cells.Item(493, 7).Value = "Zapallo italiano"
$ws.Cells.Item(493, 8).Value = "Huracán"
$ws.Cells.Item(493, 9).Value = "Primera"
$ws.Cells.Item(493, 10).Value = 120
$ws.Cells.Item(493, 11).Value = 4000
$ws.Cells.Item(493, 12).Value = 5000
$ws.Cells.Item(493, 13).Value = 4500
$ws.Cells.Item(493, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(493, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(493, 16).Value = 64
$ws.Cells.Item(493, 17).Value = 70
$ws.Cells.Item(493, 18).Value = "Hortaliza"

# Row 494 - Segunda
$ws.Cells.Item(494, 1).Value = 1
$ws.Cells.Item(494, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(494, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(494, 4).Value = 45021
$ws.Cells.Item(494, 5).Value = 15
$ws.Cells.Item(494, 6).Value = 100112032
$ws.Cells.Item(494, 7).Value = "Zapallo italiano"
$ws.Cells.Item(494, 8).Value = "Huracán"
$ws.Cells.Item(494, 9).Value = "Segunda"
$ws.Cells.Item(494, 10).Value = 160
$ws.Cells.Item(494, 11).Value = 3000
$ws.Cells.Item(494, 12).Value = 4000
$ws.Cells.Item(494, 13).Value = 3500
$ws.Cells.Item(494, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(494, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(494, 16).Value = 35
$ws.Cells.Item(494, 17).Value = 100
$ws.Cells.Item(494, 18).Value = "Hortaliza"
